$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 9.734077866766217
$ws.Cells.Item(2, 3).Value = 8.52163074024365
$ws.Cells.Item(2, 4).Value = 3.538815405798462
$ws.Cells.Item(2, 6).Value = 18.0448088309194
$ws.Cells.Item(2, 7).Value = 18.3586086075477
$ws.Cells.Item(2, 8).Value = 11.4211849205143
$ws.Cells.Item(2, 13).Value = 19.96655746346338
$ws.Cells.Item(2, 14).Value = 16.06528831776998
$ws.Cells.Item(2, 15).Value = 16.00470070490294
$ws.Cells.Item(3, 2).Value = 9.187064247943932
$ws.Cells.Item(3, 3).Value = 8.349208341236606
$ws.Cells.Item(3, 4).Value = 3.453922859612638
$ws.Cells.Item(3, 6).Value = 17.98538428828848
$ws.Cells.Item(3, 7).Value = 18.2033546309019
$ws.Cells.Item(3, 8).Value = 11.45155587626874
$ws.Cells.Item(3, 13).Value = 19.35742486105557
$ws.Cells.Item(3, 14).Value = 16.12334876121475
$ws.Cells.Item(3, 15).Value = 16.0241630408317
$ws.Cells.Item(4, 2).Value = 8.832085457688553
$ws.Cells.Item(4, 3).Value = 8.240873075138841
$ws.Cells.Item(4, 4).Value = 3.40007930463748
$ws.Cells.Item(4, 6).Value = 17.95476358268907
$ws.Cells.Item(4, 7).Value = 18.11567315398887
$ws.Cells.Item(4, 8).Value = 11.47235426441931
$ws.Cells.Item(4, 13).Value = 18.98158144228076
$ws.Cells.Item(4, 14).Value = 16.1607786420439
$ws.Cells.Item(4, 15).Value = 16.04053762878876
$ws.Cells.Item(5, 2).Value = 8.682679991230875
$ws.Cells.Item(5, 3).Value = 8.196148343104133
$ws.Cells.Item(5, 4).Value = 3.377723122017299
$ws.Cells.Item(5, 6).Value = 17.94377052313979
$ws.Cells.Item(5, 7).Value = 18.08190647181933
$ws.Cells.Item(5, 8).Value = 11.48136982733983
$ws.Cells.Item(5, 13).Value = 18.82822059646877
$ws.Cells.Item(5, 14).Value = 16.17648092523837
$ws.Cells.Item(5, 15).Value = 16.04832042757706
$ws.Cells.Item(6, 2).Value = 8.657586448005864
$ws.Cells.Item(6, 3).Value = 8.188688273510664
$ws.Cells.Item(6, 4).Value = 3.373986404726233
$ws.Cells.Item(6, 6).Value = 17.94203508354379
$ws.Cells.Item(6, 7).Value = 18.07641930000779
$ws.Cells.Item(6, 8).Value = 11.4828994507736
$ws.Cells.Item(6, 13).Value = 18.80275081044125
$ws.Cells.Item(6, 14).Value = 16.17911546046998
$ws.Cells.Item(6, 15).Value = 16.04967971179852
$ws.Cells.Item(7, 2).Value = 8.830089665780999
$ws.Cells.Item(7, 3).Value = 8.240272180466539
$ws.Cells.Item(7, 4).Value = 3.399779455211151
$ws.Cells.Item(7, 6).Value = 17.95460930157625
$ws.Cells.Item(7, 7).Value = 18.11520976005218
$ws.Cells.Item(7, 8).Value = 11.47247366579833
$ws.Cells.Item(7, 13).Value = 18.97951362803325
$ws.Cells.Item(7, 14).Value = 16.16098858728331
$ws.Cells.Item(7, 15).Value = 16.04063809975179
$ws.Cells.Item(8, 2).Value = 9.549463804829019
$ws.Cells.Item(8, 3).Value = 8.46271682341267
$ws.Cells.Item(8, 4).Value = 3.509911574311496
$ws.Cells.Item(8, 6).Value = 18.02310836680617
$ws.Cells.Item(8, 7).Value = 18.3035168737804
$ws.Cells.Item(8, 8).Value = 11.43121007725709
$ws.Cells.Item(8, 13).Value = 19.75706044040454
$ws.Cells.Item(8, 14).Value = 16.08493895702679
$ws.Cells.Item(8, 15).Value = 16.01049158582954
$ws.Cells.Item(9, 2).Value = 10.80666419555766
$ws.Cells.Item(9, 3).Value = 8.877565785374316
$ws.Cells.Item(9, 4).Value = 3.711520435113985
$ws.Cells.Item(9, 6).Value = 18.20344162701749
$ws.Cells.Item(9, 7).Value = 18.73137458162616
$ws.Cells.Item(9, 8).Value = 11.36738584069354
$ws.Cells.Item(9, 13).Value = 21.25663010709886
$ws.Cells.Item(9, 14).Value = 15.94986171888904
$ws.Cells.Item(9, 15).Value = 15.98657484562738
$ws.Cells.Item(10, 2).Value = 11.63497777497888
$ws.Cells.Item(10, 3).Value = 9.167095712670621
$ws.Cells.Item(10, 4).Value = 3.85001357949713
$ws.Cells.Item(10, 6).Value = 18.36311890449285
$ws.Cells.Item(10, 7).Value = 19.07836398459483
$ws.Cells.Item(10, 8).Value = 11.33095409397413
$ws.Cells.Item(10, 13).Value = 22.32992005094739
$ws.Cells.Item(10, 14).Value = 15.85908849856932
$ws.Cells.Item(10, 15).Value = 15.99055557176934
$ws.Cells.Item(11, 2).Value = 11.99087558848612
$ws.Cells.Item(11, 3).Value = 9.295070255519722
$ws.Cells.Item(11, 4).Value = 3.910765960267676
$ws.Cells.Item(11, 6).Value = 18.44143929348995
$ws.Cells.Item(11, 7).Value = 19.24258034688613
$ws.Cells.Item(11, 8).Value = 11.31665897593101
$ws.Cells.Item(11, 13).Value = 22.8096007841922
$ws.Cells.Item(11, 14).Value = 15.81961071694252
$ws.Cells.Item(11, 15).Value = 15.99705401819223
$ws.Cells.Item(12, 2).Value = 12.12262720835828
$ws.Cells.Item(12, 3).Value = 9.342961443855909
$ws.Cells.Item(12, 4).Value = 3.933435987930541
$ws.Cells.Item(12, 6).Value = 18.47189267872581
$ws.Cells.Item(12, 7).Value = 19.30561624133426
$ws.Cells.Item(12, 8).Value = 11.31157384716651
$ws.Cells.Item(12, 13).Value = 22.9898356147205
$ws.Cells.Item(12, 14).Value = 15.80492094006135
$ws.Cells.Item(12, 15).Value = 16.00018841979765
$ws.Cells.Item(13, 2).Value = 12.09438657421866
$ws.Cells.Item(13, 3).Value = 9.332673071014232
$ws.Cells.Item(13, 4).Value = 3.928568699395956
$ws.Cells.Item(13, 6).Value = 18.46529899997704
$ws.Cells.Item(13, 7).Value = 19.29200355322407
$ws.Cells.Item(13, 8).Value = 11.31265441896626
$ws.Cells.Item(13, 13).Value = 22.95108417139381
$ws.Cells.Item(13, 14).Value = 15.80807312344831
$ws.Cells.Item(13, 15).Value = 15.99948342729487
$ws.Cells.Item(14, 2).Value = 12.00177543948365
$ws.Cells.Item(14, 3).Value = 9.299021896740669
$ws.Cells.Item(14, 4).Value = 3.912637828569799
$ws.Cells.Item(14, 6).Value = 18.44392890746497
$ws.Cells.Item(14, 7).Value = 19.24774967751134
$ws.Cells.Item(14, 8).Value = 11.31623404025297
$ws.Cells.Item(14, 13).Value = 22.82445802530868
$ws.Cells.Item(14, 14).Value = 15.8183969856219
$ws.Cells.Item(14, 15).Value = 15.99729839283323
$ws.Cells.Item(15, 2).Value = 11.94465495830582
$ws.Cells.Item(15, 3).Value = 9.278334410811874
$ws.Cells.Item(15, 4).Value = 3.902835664543533
$ws.Cells.Item(15, 6).Value = 18.43094198092646
$ws.Cells.Item(15, 7).Value = 19.2207517260838
$ws.Cells.Item(15, 8).Value = 11.31846940856945
$ws.Cells.Item(15, 13).Value = 22.74670726457607
$ws.Cells.Item(15, 14).Value = 15.82475441426786
$ws.Cells.Item(15, 15).Value = 15.99604769071018
$ws.Cells.Item(16, 2).Value = 11.61129712185657
$ws.Cells.Item(16, 3).Value = 9.1586542577333
$ws.Cells.Item(16, 4).Value = 3.845996947791882
$ws.Cells.Item(16, 6).Value = 18.35811298768427
$ws.Cells.Item(16, 7).Value = 19.06775469544399
$ws.Cells.Item(16, 8).Value = 11.33193420008321
$ws.Cells.Item(16, 13).Value = 22.29838379692564
$ws.Cells.Item(16, 14).Value = 15.86170487080435
$ws.Cells.Item(16, 15).Value = 15.99022518228926
$ws.Cells.Item(17, 2).Value = 11.40142762464045
$ws.Cells.Item(17, 3).Value = 9.084253712122447
$ws.Cells.Item(17, 4).Value = 3.810543218221266
$ws.Cells.Item(17, 6).Value = 18.31487506516357
$ws.Cells.Item(17, 7).Value = 18.97548151772247
$ws.Cells.Item(17, 8).Value = 11.34077825441705
$ws.Cells.Item(17, 13).Value = 22.02102574349436
$ws.Cells.Item(17, 14).Value = 15.88483671291145
$ws.Cells.Item(17, 15).Value = 15.98785372359675
$ws.Cells.Item(18, 2).Value = 11.27874888659786
$ws.Cells.Item(18, 3).Value = 9.041110840915289
$ws.Cells.Item(18, 4).Value = 3.789940155594521
$ws.Cells.Item(18, 6).Value = 18.29054203583461
$ws.Cells.Item(18, 7).Value = 18.92301195916101
$ws.Cells.Item(18, 8).Value = 11.34607949429129
$ws.Cells.Item(18, 13).Value = 21.86070012087348
$ws.Cells.Item(18, 14).Value = 15.89831249792698
$ws.Cells.Item(18, 15).Value = 15.98693101199706
$ws.Cells.Item(19, 2).Value = 11.23687454908736
$ws.Cells.Item(19, 3).Value = 9.0264443997626
$ws.Cells.Item(19, 4).Value = 3.782928460540412
$ws.Cells.Item(19, 6).Value = 18.28239603845199
$ws.Cells.Item(19, 7).Value = 18.90535227962396
$ws.Cells.Item(19, 8).Value = 11.34791120352262
$ws.Cells.Item(19, 13).Value = 21.80628553147169
$ws.Cells.Item(19, 14).Value = 15.90290457432512
$ws.Cells.Item(19, 15).Value = 15.98669439376708
$ws.Cells.Item(20, 2).Value = 11.4239723125727
$ws.Cells.Item(20, 3).Value = 9.092210211026302
$ws.Cells.Item(20, 4).Value = 3.814339259409457
$ws.Cells.Item(20, 6).Value = 18.3194224666971
$ws.Cells.Item(20, 7).Value = 18.98524215783166
$ws.Cells.Item(20, 8).Value = 11.33981459855268
$ws.Cells.Item(20, 13).Value = 22.05063477646591
$ws.Cells.Item(20, 14).Value = 15.88235660554911
$ws.Cells.Item(20, 15).Value = 15.98806049922335
$ws.Cells.Item(21, 2).Value = 12.02905959370366
$ws.Cells.Item(21, 3).Value = 9.308921783339011
$ws.Cells.Item(21, 4).Value = 3.917326310802555
$ws.Cells.Item(21, 6).Value = 18.45018442468126
$ws.Cells.Item(21, 7).Value = 19.26072555306692
$ws.Cells.Item(21, 8).Value = 11.3151737093693
$ws.Cells.Item(21, 13).Value = 22.86169078901029
$ws.Cells.Item(21, 14).Value = 15.81535758302808
$ws.Cells.Item(21, 15).Value = 15.997921916728
$ws.Cells.Item(22, 2).Value = 12.40691822969618
$ws.Cells.Item(22, 3).Value = 9.44721930366029
$ws.Cells.Item(22, 4).Value = 3.982673392376071
$ws.Cells.Item(22, 6).Value = 18.54026836032912
$ws.Cells.Item(22, 7).Value = 19.44569721152682
$ws.Cells.Item(22, 8).Value = 11.30098222734927
$ws.Cells.Item(22, 13).Value = 23.38347491266276
$ws.Cells.Item(22, 14).Value = 15.77308234326879
$ws.Cells.Item(22, 15).Value = 16.00829238441817
$ws.Cells.Item(23, 2).Value = 12.20686155514831
$ws.Cells.Item(23, 3).Value = 9.373722864276909
$ws.Cells.Item(23, 4).Value = 3.94797959574791
$ws.Cells.Item(23, 6).Value = 18.49177350156346
$ws.Cells.Item(23, 7).Value = 19.34654545121056
$ws.Cells.Item(23, 8).Value = 11.30838129995136
$ws.Cells.Item(23, 13).Value = 23.10580154587252
$ws.Cells.Item(23, 14).Value = 15.79550751774723
$ws.Cells.Item(23, 15).Value = 16.00239861941952
$ws.Cells.Item(24, 2).Value = 11.41378615611549
$ws.Cells.Item(24, 3).Value = 9.088614224713082
$ws.Cells.Item(24, 4).Value = 3.812623754253301
$ws.Cells.Item(24, 6).Value = 18.31736494947204
$ws.Cells.Item(24, 7).Value = 18.9808275624408
$ws.Cells.Item(24, 8).Value = 11.34024959251779
$ws.Cells.Item(24, 13).Value = 22.0372512249982
$ws.Cells.Item(24, 14).Value = 15.88347731072192
$ws.Cells.Item(24, 15).Value = 15.98796564333937
$ws.Cells.Item(25, 2).Value = 10.48320681756999
$ws.Cells.Item(25, 3).Value = 8.767865137052091
$ws.Cells.Item(25, 4).Value = 3.658612256634702
$ws.Cells.Item(25, 6).Value = 18.14981230815581
$ws.Cells.Item(25, 7).Value = 18.6096583720293
$ws.Cells.Item(25, 8).Value = 11.38281773521265
$ws.Cells.Item(25, 13).Value = 20.85507299108772
$ws.Cells.Item(25, 14).Value = 15.98490938150241
$ws.Cells.Item(25, 15).Value = 15.98926452440247

Write-Output "Updated 216 cells in loading_percent sheet (Case_3_200, 380 kV)"
